# Insert a new weekly price record as row 119, pushing the existing
# rows 119-138 down to 120-139 (data for "Femacal de La Calera - Achicoria").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 119; Excel shifts rows 119:138 -> 120:139
# and inherits formatting (including the date style on column D) from the
# surrounding rows.
$ws.Rows.Item(119).Insert()

# Populate the new row 119 with the new weekly record.
$ws.Range("A119").Value = 3
$ws.Range("B119").Value = "Femacal de La Calera"
$ws.Range("C119").Value = "Coquimbo"
$ws.Range("D119").Value = 44504
$ws.Range("E119").Value = 5
$ws.Range("F119").Value = 100112010
$ws.Range("G119").Value = "Achicoria"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 125
$ws.Range("K119").Value = 6000
$ws.Range("L119").Value = 6500
$ws.Range("M119").Value = 6240
$ws.Range("N119").Value = "$/caja 16 unidades"
$ws.Range("O119").Value = "Provincia de Quillota"
$ws.Range("P119").Value = 390
$ws.Range("Q119").Value = 16
$ws.Range("R119").Value = "Hortaliza"
